$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "last_update"
$ws.Range("D1").Value = "location"

$ws.Range("C14").Select()
